$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "in the hope that you and your generation might discover a future where technology and personal data drive human flourishing more than corporate profit." "in the hope that you and your generation might experience a future where technology empowers individuals and personal data drives human flourishing more than corporate profit."

Replace-Text "And also to Rob Comber, Simon Bowen" "And also, to Rob Comber, Simon Bowen"

Replace-Text "Finally I would also like to thank, in no particular order:" "Finally, I would also like to thank, in no particular order:"

Replace-Text "- Home Interviewing: Card Sorting With a Family in Their Living Room" "- Home Interviewing: Card Sorting with a Family in Their Living Room"

Replace-Text "Created through Discussion From a Palette of Possible Parent and Staff Actions" "Created through Discussion from a Palette of Possible Parent and Staff Actions"

Replace-Text "- Mockup of a Unified TV Viewing History Interface" "- Mock-up of a Unified TV Viewing History Interface"

Replace-Text "- Mockup of a Unified Interface for a Vacation" "- Mock-up of a Unified Interface for a Vacation"

Replace-Text "- Mockup of Life Information Presented in a PDS Interface" "- Mock-up of Life Information Presented in a PDS Interface"

Replace-Text "- Mockup: Browsing By Areas of Life" "- Mock-up: Browsing by Areas of Life"

Replace-Text "- Determining The Nature of a Piece of Data" "- Determining the Nature of a Piece of Data"

Replace-Text "- Example Taxonomies For Life Information Navigation" "- Example Taxonomies for Life Information Navigation"

Replace-Text "- Screenshot From Quirkos During Coding Process" "- Screenshot from Quirkos During Coding Process"

Replace-Text "- Screenshot From Quirkos at End of Coding Process" "- Screenshot from Quirkos at End of Coding Process"

Replace-Text "- Screenshot From Workflowy During Theme Construction" "- Screenshot from Workflowy During Theme Construction"
